# Edit "traceability matrix.xlsx" - Sheet1 (Table2, range A1:H9)
# Updates the "Filter flight by altitude" (#2) test rows and the
# "Info Boxes" (#3) section per the author's revision.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: SAATM-2-1 (altitude 0-16000 filtering test) ---
$rsquo = [char]0x2019
$endash = [char]0x2013
$ws.Range("E5").Value = "Flight" + $rsquo + "s altitude 0 " + $endash + " 16000 ft. filtering test. "
$ws.Range("F5").Value = "Test whether our software can filter flight altitude from 0 ft. to 16000 ft."

# --- Row 6: SAATM-2-2 (minimum altitude filtering test) ---
$ws.Range("E6").Value = "Flight's minimum altitude filtering test (10000 ft)  "
$ws.Range("F6").Value = "Test whether our software can filter flights by minimum altitude of 10000 ft."

# --- Row 7: SAATM-2-3 (maximum altitude filtering test) ---
$ws.Range("E7").Value = "Flight's maximum altitude filtering test (16000 ft)  "
$ws.Range("F7").Value = "Test whether our software can filter flights by maximum altitude of 16000 ft."

# --- Row 8: now SAATM-2-4 (clear all altitude filters test); drop the old
#     "#3 / Info Boxes" section header cells that used to live here ---
$ws.Range("A8:D8").ClearContents()
$ws.Range("E8").Value = "Clear all altitude filters test"
$ws.Range("F8").Value = "Test whether our software can clear all altitude filters and show all flights."
$ws.Range("G8").Value = "SAATM-2-4"
$ws.Range("H8").Value = "Pass"

# --- Row 9: the "#3 / Info Boxes" section header now starts here ---
$ws.Range("A9").Value = "#3"
$ws.Range("B9").Value = "Flight info boxes"
$ws.Range("C9").Value = "Allow user to hide or show info boxes."
$ws.Range("D9").Value = "Info Boxes"
$ws.Range("E9").Value = "Show info boxes."
$ws.Range("F9").Value = "Test whether our software can show the info boxes."
$ws.Range("G9").Value = "SAATM-3-1"
$ws.Range("H9").Value = "Pass"

# Match the author's final cursor/selection position.
$ws.Range("F14").Select()
